# The workbook originally has two sheets:
#   Worksheets.Item(1) = "hotel_info"  (9 header cols + 1 data row)
#   Worksheets.Item(2) = "review_info" (24 header cols, no data row)
#
# The edit swaps the sheet names AND restructures their contents:
#   Worksheets.Item(1) becomes "review_info" with the full 25-column
#     review header row (STR + 24 review fields) and no data row.
#   Worksheets.Item(2) becomes "hotel_info" with a 10-column header row
#     (STR, Hotel_Name, State, City, Zip, TA_ReviewURL,
#      Tripadvisor_Hotel_Name, English_Reviews_num, Local_Rank,
#      Total_Reviews_num) plus one data row, now including a new
#      "State" column ("Louisiana").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Wipe both sheets completely so no stale cells remain ---
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# --- Rename to final names (via a temp name to avoid a collision,
#     since the two sheets are swapping names) ---
$ws1.Name = "__tmp_swap__"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"

# --- review_info header row (A1:Y1) ---
$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- hotel_info header row (A1:J1) ---
$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# --- hotel_info data row (A2:J2) ---
$ws2.Cells.Item(2, 1).Value = 53513
$ws2.Cells.Item(2, 2).Value = "Marriott New Orleans Downtown @ The Convention Center"
$ws2.Cells.Item(2, 3).Value = "Louisiana"
$ws2.Cells.Item(2, 4).Value = "New Orleans"
$ws2.Cells.Item(2, 5).Value = 70130
$ws2.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d553884-Reviews-New_Orleans_Downtown_Marriott_at_the_Convention_Center-New_Orleans_Louisiana.html"
$ws2.Cells.Item(2, 7).Value = "New Orleans Downtown Marriott at the Convention Center"

# English_Reviews_num / Local_Rank / Total_Reviews_num are stored as TEXT
# (not numbers) in the source data, even though they look numeric. Force
# the cells to text via a temporary "@" number format (otherwise Excel
# auto-detects the numeric-looking string and stores it as a number),
# then clear the format again so no extra formatting lingers on the cell.
$ws2.Cells.Item(2, 8).NumberFormat = "@"
$ws2.Cells.Item(2, 8).Value = "866"
$ws2.Cells.Item(2, 8).ClearFormats()

$ws2.Cells.Item(2, 9).NumberFormat = "@"
$ws2.Cells.Item(2, 9).Value = "62"
$ws2.Cells.Item(2, 9).ClearFormats()

$ws2.Cells.Item(2, 10).NumberFormat = "@"
$ws2.Cells.Item(2, 10).Value = "900"
$ws2.Cells.Item(2, 10).ClearFormats()

$wb.Worksheets.Item(1).Select()
